# team quest 2 - game play
#
# 1. Delete the extra "튜토리얼" rectangle (id=6, "직사각형 5") duplicated on
#    slide 1.
# 2. Bump every cached datetimeFigureOut field ("2024-07-30" -> "2024-07-31")
#    across the slide master and all 11 slide layouts.

$p = $ppt.ActivePresentation

# --- 1. Remove the duplicate "튜토리얼" rectangle on slide 1 -----------------
$s1 = $p.Slides.Item(1)
for ($i = $s1.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s1.Shapes.Item($i)
    if ($sh.Name -eq "직사각형 5" -and $sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "튜토리얼") {
        $sh.Delete()
    }
}

# --- helper: refresh a cached date field's text on a shape collection -------
function Update-DateField($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "2024-07-30") {
            $sh.TextFrame.TextRange.Text = "2024-07-31"
        }
    }
}

# --- 2. Slide master ----------------------------------------------------
$master = $p.SlideMaster
Update-DateField $master.Shapes

# --- 3. Every slide layout ------------------------------------------------
$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DateField $layouts.Item($L).Shapes
}
